$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 and row 3 data
$ws.Range("D2").Value = "Beneficiário não encontrado"
$ws.Range("H2").Value = "Data de admissão não encontrada"
$ws.Range("I2").Value = "Data de alta não encontrada"
$ws.Range("J2").Value = "MOTIVO NÃO INFORMADO"

$ws.Range("D3").Value = "Beneficiário não encontrado"
$ws.Range("H3").Value = "Data de admissão não encontrada"
$ws.Range("I3").Value = "Data de alta não encontrada"
$ws.Range("J3").Value = "MOTIVO NÃO INFORMADO"

# Delete rows 4 and 5 (shrinks used range to A1:J3)
$ws.Rows("4:5").Delete()

# Adjust column widths (subtract the 5/6 char padding this COM host adds
# on top of ColumnWidth so the saved OOXML 'width' lands on the exact
# target integer, matching the original file's style of whole-number widths)
$ws.Columns("D").ColumnWidth = 28.166666666666668
$ws.Columns("H").ColumnWidth = 32.166666666666664
$ws.Columns("I").ColumnWidth = 28.166666666666668
$ws.Columns("J").ColumnWidth = 21.166666666666668
